# QB Website updated 12/13
# - Append "83" to the existing Occurrence list shared by E4/E7 (making it
#   "...69, 75, 83"), and append "83, 87" to the Occurrence list shown in
#   E2/E3/E5/E6 (making it "...69, 75, 83, 87").
# - Update the sheet view: scroll so column A (not C) is the left-most
#   visible column, keep row 6 as the top visible row, and set zoom to 83%.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: extend the list all six Occurrence cells currently share -----
# All of E2:E7 start out pointing at the same shared string
# ("21, 23, ..., 69, 75"). Rewrite every one of them to the same new text so
# the underlying shared string is updated in place (keeping its identity)
# instead of forking a new, disconnected entry.
$extended = "21, 23, 25, 27, 31, 33, 35, 37, 39, 43, 47, 51, 55, 59, 63, 67, 69, 75, 83"
$ws.Range("E2").Value = $extended
$ws.Range("E3").Value = $extended
$ws.Range("E4").Value = $extended
$ws.Range("E5").Value = $extended
$ws.Range("E6").Value = $extended
$ws.Range("E7").Value = $extended

# --- Step 2: give E2, E3, E5, E6 one more entry (87), leaving E4 and E7 ----
# pointing at the "...75, 83" text set above.
$furtherExtended = "21, 23, 25, 27, 31, 33, 35, 37, 39, 43, 47, 51, 55, 59, 63, 67, 69, 75, 83, 87"
$ws.Range("E2").Value = $furtherExtended
$ws.Range("E3").Value = $furtherExtended
$ws.Range("E5").Value = $furtherExtended
$ws.Range("E6").Value = $furtherExtended

# --- Step 3: update the view ----------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 6
$win.Zoom = 83
